# Auto-generated edit script: updates Betfair odds values for rows 2-22 (columns F:AO)
# as described in the commit "Atualizando o arquivo XLSX".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.09
$ws.Range("G2").Value = 7
$ws.Range("H2").Value = 1.75
$ws.Range("I2").Value = 2.08
$ws.Range("J2").Value = 2.72
$ws.Range("K2").Value = 3.7
$ws.Range("L2").Value = 1.52
$ws.Range("M2").Value = 1.12
$ws.Range("N2").Value = 2.12
$ws.Range("O2").Value = 1.56
$ws.Range("P2").Value = 1.41
$ws.Range("Q2").Value = 2.3
$ws.Range("R2").Value = 1.14
$ws.Range("T2").Value = 2.3
$ws.Range("U2").Value = 1.53
$ws.Range("V2").Value = 1.92
$ws.Range("W2").Value = 1.16
$ws.Range("X2").Value = 9.2
$ws.Range("Y2").Value = 6.2
$ws.Range("Z2").Value = 10.5
$ws.Range("AA2").Value = 25
$ws.Range("AB2").Value = 15
$ws.Range("AC2").Value = 8.8
$ws.Range("AD2").Value = 13
$ws.Range("AE2").Value = 34
$ws.Range("AF2").Value = 55
$ws.Range("AG2").Value = 30
$ws.Range("AH2").Value = 36
$ws.Range("AI2").Value = 90
$ws.Range("AJ2").Value = 270
$ws.Range("AK2").Value = 170
$ws.Range("AL2").Value = 210
$ws.Range("AN2").Value = 350
$ws.Range("AO2").Value = 29

# Row 3
$ws.Range("F3").Value = 1.56
$ws.Range("J3").Value = 4.3
$ws.Range("K3").Value = 4.7
$ws.Range("N3").Value = 4
$ws.Range("P3").Value = 2.02
$ws.Range("S3").Value = 3.1
$ws.Range("T3").Value = 1.93
$ws.Range("U3").Value = 1.9
$ws.Range("W3").Value = 2.58
$ws.Range("X3").Value = 20
$ws.Range("Z3").Value = 65
$ws.Range("AA3").Value = 250
$ws.Range("AC3").Value = 10.5
$ws.Range("AJ3").Value = 17.5
$ws.Range("AN3").Value = 10

# Row 4
$ws.Range("F4").Value = 5.5
$ws.Range("G4").Value = 5.7
$ws.Range("H4").Value = 1.6
$ws.Range("I4").Value = 1.62
$ws.Range("Q4").Value = 1.67
$ws.Range("S4").Value = 2.68
$ws.Range("T4").Value = 1.77
$ws.Range("V4").Value = 2.62
$ws.Range("W4").Value = 1.21
$ws.Range("X4").Value = 25
$ws.Range("AA4").Value = 15.5
$ws.Range("AB4").Value = 24
$ws.Range("AD4").Value = 9.8
$ws.Range("AE4").Value = 15
$ws.Range("AF4").Value = 48
$ws.Range("AG4").Value = 22
$ws.Range("AK4").Value = 70
$ws.Range("AN4").Value = 65
$ws.Range("AO4").Value = 7

# Row 5
$ws.Range("G5").Value = 2.26
$ws.Range("J5").Value = 3.7
$ws.Range("W5").Value = 1.79
$ws.Range("AN5").Value = 13.5

# Row 6
$ws.Range("G6").Value = 1.7
$ws.Range("K6").Value = 5.2

# Row 7
$ws.Range("F7").Value = 8.2
$ws.Range("G7").Value = 8.4
$ws.Range("H7").Value = 1.38
$ws.Range("I7").Value = 1.39
$ws.Range("U7").Value = 2.5
$ws.Range("V7").Value = 3.55
$ws.Range("AJ7").Value = 250

# Row 8
$ws.Range("F8").Value = 2.88
$ws.Range("H8").Value = 2.6
$ws.Range("K8").Value = 3.7
$ws.Range("P8").Value = 1.76
$ws.Range("R8").Value = 1.28
$ws.Range("S8").Value = 3.85
$ws.Range("T8").Value = 1.82
$ws.Range("U8").Value = 2
$ws.Range("V8").Value = 1.56
$ws.Range("AL8").Value = 50
$ws.Range("AM8").Value = 130

# Row 9
$ws.Range("F9").Value = 1.9
$ws.Range("G9").Value = 2.08
$ws.Range("I9").Value = 4.9
$ws.Range("J9").Value = 3.5
$ws.Range("K9").Value = 4.5
$ws.Range("M9").Value = 1.01
$ws.Range("P9").Value = 1.83
$ws.Range("R9").Value = 1.32
$ws.Range("S9").Value = 2.26
$ws.Range("T9").Value = 1.81
$ws.Range("W9").Value = 1.92
$ws.Range("AD9").Value = 20
$ws.Range("AE9").Value = 65
$ws.Range("AJ9").Value = 28
$ws.Range("AK9").Value = 25
$ws.Range("AL9").Value = 44

# Row 10
$ws.Range("F10").Value = 5.5
$ws.Range("I10").Value = 1.65
$ws.Range("J10").Value = 3.75
$ws.Range("S10").Value = 2.66
$ws.Range("V10").Value = 2.5

# Row 11
$ws.Range("G11").Value = 1.68
$ws.Range("H11").Value = 5.2
$ws.Range("P11").Value = 2.6
$ws.Range("Q11").Value = 1.5
$ws.Range("R11").Value = 1.65
$ws.Range("S11").Value = 2.24
$ws.Range("AA11").Value = 150
$ws.Range("AN11").Value = 7.2
$ws.Range("AO11").Value = 60

# Row 12
$ws.Range("F12").Value = 1.79
$ws.Range("G12").Value = 1.87
$ws.Range("H12").Value = 4.6
$ws.Range("J12").Value = 3.85
$ws.Range("K12").Value = 4.2
$ws.Range("V12").Value = 1.25
$ws.Range("W12").Value = 2.16
$ws.Range("AG12").Value = 10.5
$ws.Range("AJ12").Value = 20

# Row 13
$ws.Range("G13").Value = 3.25
$ws.Range("H13").Value = 2.34
$ws.Range("I13").Value = 2.48
$ws.Range("K13").Value = 3.9
$ws.Range("V13").Value = 1.67

# Row 14
$ws.Range("G14").Value = 1.46
$ws.Range("H14").Value = 8.6
$ws.Range("K14").Value = 5.8
$ws.Range("L14").Value = 1.29
$ws.Range("P14").Value = 2.16
$ws.Range("Q14").Value = 1.7
$ws.Range("S14").Value = 2.74
$ws.Range("T14").Value = 1.98
$ws.Range("U14").Value = 1.83
$ws.Range("W14").Value = 3.15
$ws.Range("Z14").Value = 100
$ws.Range("AA14").Value = 380
$ws.Range("AJ14").Value = 13.5
$ws.Range("AL14").Value = 44
$ws.Range("AO14").Value = 230

# Row 15
$ws.Range("G15").Value = 2.74
$ws.Range("I15").Value = 3.55
$ws.Range("K15").Value = 3.5

# Row 16
$ws.Range("F16").Value = 3.05
$ws.Range("G16").Value = 3.1
$ws.Range("H16").Value = 2.42
$ws.Range("I16").Value = 2.44
$ws.Range("N16").Value = 4.5
$ws.Range("O16").Value = 1.26
$ws.Range("S16").Value = 3
$ws.Range("V16").Value = 1.68
$ws.Range("W16").Value = 1.48
$ws.Range("Z16").Value = 17.5
$ws.Range("AD16").Value = 11
$ws.Range("AF16").Value = 22
$ws.Range("AG16").Value = 13
$ws.Range("AN16").Value = 24
$ws.Range("AO16").Value = 17

# Row 17
$ws.Range("F17").Value = 4.4
$ws.Range("G17").Value = 4.6
$ws.Range("H17").Value = 1.95
$ws.Range("J17").Value = 3.65
$ws.Range("Q17").Value = 2.06
$ws.Range("S17").Value = 3.8
$ws.Range("AF17").Value = 32
$ws.Range("AG17").Value = 17.5

# Row 18
$ws.Range("I18").Value = 2.68
$ws.Range("P18").Value = 1.77
$ws.Range("S18").Value = 4.2
$ws.Range("U18").Value = 2.06
$ws.Range("V18").Value = 1.59
$ws.Range("X18").Value = 11.5
$ws.Range("Y18").Value = 9.6
$ws.Range("AG18").Value = 13.5

# Row 19
$ws.Range("F19").Value = 4.6
$ws.Range("P19").Value = 2.5
$ws.Range("R19").Value = 1.6
$ws.Range("S19").Value = 2.6
$ws.Range("U19").Value = 2.48
$ws.Range("V19").Value = 2.22
$ws.Range("Z19").Value = 12.5

# Row 20
$ws.Range("F20").Value = 1.5
$ws.Range("G20").Value = 1.51
$ws.Range("H20").Value = 7
$ws.Range("I20").Value = 7.2
$ws.Range("J20").Value = 5
$ws.Range("K20").Value = 5.1
$ws.Range("L20").Value = 1.25
$ws.Range("O20").Value = 1.14
$ws.Range("P20").Value = 3.05
$ws.Range("Q20").Value = 1.46
$ws.Range("R20").Value = 1.83
$ws.Range("S20").Value = 2.14
$ws.Range("V20").Value = 1.16
$ws.Range("W20").Value = 2.96
$ws.Range("X20").Value = 32
$ws.Range("Y20").Value = 40
$ws.Range("Z20").Value = 70
$ws.Range("AA20").Value = 190
$ws.Range("AC20").Value = 12.5
$ws.Range("AD20").Value = 25
$ws.Range("AE20").Value = 75
$ws.Range("AF20").Value = 11.5
$ws.Range("AG20").Value = 10.5
$ws.Range("AH20").Value = 18.5
$ws.Range("AI20").Value = 70
$ws.Range("AJ20").Value = 15
$ws.Range("AM20").Value = 75
$ws.Range("AN20").Value = 4.8
$ws.Range("AO20").Value = 980

# Row 21
$ws.Range("N21").Value = 8.6
$ws.Range("P21").Value = 3.5
$ws.Range("Q21").Value = 1.39
$ws.Range("R21").Value = 1.99
$ws.Range("S21").Value = 1.98
$ws.Range("T21").Value = 1.72
$ws.Range("AB21").Value = 15
$ws.Range("AG21").Value = 11
$ws.Range("AH21").Value = 23
$ws.Range("AN21").Value = 3.6
$ws.Range("AO21").Value = 90

# Row 22
$ws.Range("F22").Value = 2.58
$ws.Range("G22").Value = 2.62
$ws.Range("H22").Value = 2.8
$ws.Range("I22").Value = 2.84
$ws.Range("Q22").Value = 1.8
$ws.Range("U22").Value = 2.4
$ws.Range("V22").Value = 1.54
$ws.Range("W22").Value = 1.61
$ws.Range("X22").Value = 16.5
